# Append: 2026-02-01 12:43 JST
# Scraper run added two new postings near the top, refreshed the
# "taken at" timestamp on every existing row, and appended one more
# posting at the bottom of the "ランサーズ" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-01 12:43:33"

# Drop the existing hyperlinks up front - row numbers are about to shift
# and we will re-create them at their final locations below.
$ws.Hyperlinks.Delete()

# Make room for the two new postings: push the current rows 3-5 down to
# rows 5-7 (row 2 keeps its place).
$ws.Rows.Item(3).Resize(2).Insert()

# --- Row 2 : existing posting, only the scrape timestamp changes ---
$ws.Range("A2").Value = $newTimestamp

# --- Row 3 : new posting ---
$url3 = "https://www.lancers.jp/work/detail/5483343"
$ws.Range("A3").Value = $newTimestamp
$ws.Range("B3").Value = "【急募】ウェブサイトのAIチャットサポートの実装とLINE公式アカウントのAIチャットサポート"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = $url3
$ws.Range("G3").Value = 338
$ws.Range("H3").Value = "🔥AI,Ai ◇サイト"

# --- Row 4 : new posting ---
$url4 = "https://www.lancers.jp/work/detail/5483345"
$ws.Range("A4").Value = $newTimestamp
$ws.Range("B4").Value = "【急募】人事評価・賞与計算を自動化する社内向けWebシステム開発"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = $url4
$ws.Range("G4").Value = 203
$ws.Range("H4").Value = "◆開発,システム開発"

# --- Row 5 : previously row 3, timestamp refreshed, content unchanged ---
$url5 = "https://www.lancers.jp/work/detail/5483207"
$ws.Range("A5").Value = $newTimestamp
$ws.Range("B5").Value = "本人顔ベースのリアルタイム顔変換システム開発"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = $url5
$ws.Range("G5").Value = 125
$ws.Range("H5").Value = "◆開発,システム開発"

# --- Row 6 : previously row 4, timestamp refreshed, content unchanged ---
$url6 = "https://www.lancers.jp/work/detail/5483306"
$ws.Range("A6").Value = $newTimestamp
$ws.Range("B6").Value = "【PM/フルスタックエンジニア】新規SaaS開発のパートナー募集"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = $url6
$ws.Range("G6").Value = 75
$ws.Range("H6").Value = "◆開発"

# --- Row 7 : previously row 5, timestamp refreshed, content unchanged ---
$url7 = "https://www.lancers.jp/work/detail/5483311"
$ws.Range("A7").Value = $newTimestamp
$ws.Range("B7").Value = "カフェ掲載プラットフォーム「チャヤドコ」開発(要件定義~ベータ版リリース)"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = $url7
$ws.Range("G7").Value = 68
$ws.Range("H7").Value = "◆開発"

# --- Row 8 : brand new posting appended at the end ---
$url8 = "https://www.lancers.jp/work/detail/5483295"
$ws.Range("A8").Value = $newTimestamp
$ws.Range("B8").Value = "【月次継続】hacomono・Stripeデータ管理のプロを求む!"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = $url8
$ws.Range("G8").Value = 38
$ws.Range("H8").Value = "◇管理"

# Re-create the URL hyperlinks for every data row at their final position.
$url2 = "https://www.lancers.jp/work/detail/5483313"
$ws.Hyperlinks.Add($ws.Range("F2"), $url2)
$ws.Hyperlinks.Add($ws.Range("F3"), $url3)
$ws.Hyperlinks.Add($ws.Range("F4"), $url4)
$ws.Hyperlinks.Add($ws.Range("F5"), $url5)
$ws.Hyperlinks.Add($ws.Range("F6"), $url6)
$ws.Hyperlinks.Add($ws.Range("F7"), $url7)
$ws.Hyperlinks.Add($ws.Range("F8"), $url8)

# Column B and H got a bit wider to fit the new, longer titles/tags.
# ColumnWidth goes through Excel's character->pixel->character rounding,
# so feed it the pre-image of the integer widths we actually want.
$ws.Columns.Item(2).ColumnWidth = 48.166666666666664
$ws.Columns.Item(8).ColumnWidth = 12.166666666666666
